# Atualizado por script em 21-12-2023 20:46
#
# This script:
#  1) Swaps the F:V (home..url_partida) payload between row pairs
#     114/115, 163/164 and 214/215 (the A:E "match identity" columns -
#     Indice/pais/torneio/temporada/data_partida - stay put on their
#     original row).
#  2) Appends two new match rows (229 and 230) at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowPayload {
    param($ws, [int]$rowA, [int]$rowB)

    # Columns F..V (6..22) hold the match payload that needs to trade places.
    for ($col = 6; $col -le 22; $col++) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)
        $tmp = $cellA.Value2
        $cellA.Value2 = $cellB.Value2
        $cellB.Value2 = $tmp
    }
}

Swap-RowPayload $ws 114 115
Swap-RowPayload $ws 163 164
Swap-RowPayload $ws 214 215

# --- Append the two new rows at the bottom (229 and 230) -----------------
# Clone formatting from the last existing data row (228) so the new rows
# keep the same styles (bold/bordered index column, datetime-formatted
# date column, etc.)
$ws.Range("A228:V228").Copy()
$ws.Range("A229:V230").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A229").Value = 228
$ws.Range("B229").Value = "spain"
$ws.Range("C229").Value = "laliga2"
$ws.Range("D229").Value = "2023-2024"
$ws.Range("E229").Value = 45281.79166666666
$ws.Range("F229").Value = "Huesca"
$ws.Range("G229").Value = 3
$ws.Range("H229").Value = "FC Cartagena SAD"
$ws.Range("I229").Value = 0
$ws.Range("J229").Value = 2.16
$ws.Range("K229").Value = "17/12/2023 16:42"
$ws.Range("L229").Value = 2.11
$ws.Range("M229").Value = "21/12/2023 18:59"
$ws.Range("N229").Value = 3.11
$ws.Range("O229").Value = "17/12/2023 16:42"
$ws.Range("P229").Value = 2.87
$ws.Range("Q229").Value = "21/12/2023 18:59"
$ws.Range("R229").Value = 3.92
$ws.Range("S229").Value = "17/12/2023 16:42"
$ws.Range("T229").Value = 4.81
$ws.Range("U229").Value = "21/12/2023 18:59"
$ws.Range("V229").Value = "https://www.betexplorer.com/football/spain/laliga2/huesca-fc-cartagena-sad/SAf5sFng/"

$ws.Range("A230").Value = 229
$ws.Range("B230").Value = "spain"
$ws.Range("C230").Value = "laliga2"
$ws.Range("D230").Value = "2023-2024"
$ws.Range("E230").Value = 45281.79166666666
$ws.Range("F230").Value = "Villarreal B"
$ws.Range("G230").Value = 1
$ws.Range("H230").Value = "R. Oviedo"
$ws.Range("I230").Value = 1
$ws.Range("J230").Value = 3.09
$ws.Range("K230").Value = "18/12/2023 20:42"
$ws.Range("L230").Value = 2.99
$ws.Range("M230").Value = "21/12/2023 18:56"
$ws.Range("N230").Value = 3.04
$ws.Range("O230").Value = "18/12/2023 20:42"
$ws.Range("P230").Value = 3.05
$ws.Range("Q230").Value = "21/12/2023 18:56"
$ws.Range("R230").Value = 2.6
$ws.Range("S230").Value = "18/12/2023 20:42"
$ws.Range("T230").Value = 2.74
$ws.Range("U230").Value = "21/12/2023 18:56"
$ws.Range("V230").Value = "https://www.betexplorer.com/football/spain/laliga2/villarreal-r-oviedo/IamIvDWB/"
